$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 (shifts Lx row from row2 -> row3)
$ws.Rows.Item(2).Insert()

# Insert two new rows before what is now row 4 (shifts Ly row from row4 -> row6)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Insert a new row after the Ly row (row6), i.e. before row7
$ws.Rows.Item(7).Insert()

# Clear any formatting inherited from the row above during insert
$ws.Range("A2:D2").ClearFormats()
$ws.Range("A4:D4").ClearFormats()
$ws.Range("A5:D5").ClearFormats()
$ws.Range("A7:D7").ClearFormats()

# Fill the new rows with -99 values
$ws.Range("A2:D2").Value = -99
$ws.Range("A4:D4").Value = -99
$ws.Range("A5:D5").Value = -99
$ws.Range("A7:D7").Value = -99

$ws.Range("A7:D7").Select()
